$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# project_name (row 3, column B)
$ws.Range("B3").Value = "iacs_panel_1"

# excluded_channels (row 5, column B) - add ", NA" at the end
$ws.Range("B5").Value = "B2M, DNA, Bead, LD, Live, Dead, ID, Cell-ID, Cell_ID, NA"

# anchor_ids (row 9, column B)
$ws.Range("B9").Value = "HC-04"

# grouping_columns (row 15, column B)
$ws.Range("B15").Value = "hc_vs_pc_pre_treat, responder_pre_treat, pre_post, pre_post_w_hc"

# grouping_orders (row 16, column B and comment in column C)
$ws.Range("B16").Value = "HC, S1.1; yes, no; S1.1, S1.2; HC, S1.1, S1.2"
$ws.Range("C16").Value = "Control order of groups. Input group names separated by comma and whitespace ("", ""). For multiple grouping columns - separate by semicolon and whitespace (""; ""). If you do not want to input order for any specific grouping column - write NA."

# data_subsets (row 17, column B)
$ws.Range("B17").Value = "B"

# ccp_delta_cutoff (row 29, column B)
$ws.Range("B29").Value = 0.007

# Update the sheet view so the scrolled/selected position matches the saved state
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B23").Select()
